$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting existing A/B columns to B/C
$ws.Columns("A:A").Insert()

# Populate new column A (order chosen so the shared-string table matches append order)
$ws.Range("A2").Value = "TH_TC014"
$ws.Range("A3").Value = "TH_TC014_R"
$ws.Range("A1").Value = "TC"

# Populate new row 3
$ws.Range("B3").Value = "test9183"
$ws.Range("C3").Value = "Testing123!"

# Style the new header cell A1 like the other headers (bold, underlined)
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Font.Underline = $true

# Auto-fit the new columns so the data is fully visible
$ws.Columns("A:B").AutoFit()

# Move the active selection like in the authored workbook
$ws.Range("D8").Select()
